# Chapter 1: revisions (v0.3.2)
# Add a new day's row (2020-10-31, serial 44142) to the "writing" log table,
# extend the two chart series that plot it, and leave the selection on the
# newly-added row like the author did.

$wb = $excel.ActiveWorkbook

$dataSheet  = $wb.Worksheets.Item("writing")
$chartSheet = $wb.Worksheets.Item("dashboard")

# Remember which sheet was active so we can restore it at the end - adding
# data to "writing" / selecting a cell there should not steal the active
# tab away from "dashboard".
$originalActiveSheet = $wb.ActiveSheet.Name

# --- 1. Grow the Excel Table (ListObject) by one row -----------------------
# This keeps the table's `ref`/`autoFilter` ranges (and the sheet dimension)
# in sync automatically, the same way typing into the row below an existing
# Table does in the real app.
$table = $dataSheet.ListObjects.Item("Table1")
$table.ListRows.Add() | Out-Null

# --- 2. Fill in the new row's data ------------------------------------------
$newRow = 12
$prevRow = $newRow - 1

# Carry the date column's number format (m/d/yyyy, style index shared with
# the rows above) down onto the new row before writing its value, same as
# Excel does when you extend a table by typing into the row underneath it.
$dataSheet.Range("A" + $prevRow).Copy() | Out-Null
$dataSheet.Range("A" + $newRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$dataSheet.Range("A" + $newRow).Value = 44142
$dataSheet.Range("B" + $newRow).Value = 239
$dataSheet.Range("C" + $newRow).Value = 87
$dataSheet.Range("D" + $newRow).Value = 547
$dataSheet.Range("E" + $newRow).Value = 7413
$dataSheet.Range("F" + $newRow).Value = 297
$dataSheet.Range("G" + $newRow).Value = 6
$dataSheet.Range("H" + $newRow).Value = 5
$dataSheet.Range("I" + $newRow).Value = 5
$dataSheet.Range("J" + $newRow).Formula = "=SUM(B" + $newRow + ":I" + $newRow + ")"
$dataSheet.Range("K" + $newRow).Formula = "=J" + $newRow + "-J" + ($newRow - 1)

# --- 3. Extend the dashboard chart's two series to include the new row -----
$chartObj = $chartSheet.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection()

$seriesK = $series.Item(1)
$seriesK.Formula = "=SERIES(writing!`$K`$1,writing!`$A`$2:`$A`$" + $newRow + ",writing!`$K`$2:`$K`$" + $newRow + ",1)"

$seriesJ = $series.Item(2)
$seriesJ.Formula = "=SERIES(writing!`$J`$1,writing!`$A`$2:`$A`$" + $newRow + ",writing!`$J`$2:`$J`$" + $newRow + ",2)"

# --- 4. Match the author's final selection (cell H12 on "writing") ---------
$dataSheet.Activate() | Out-Null
$dataSheet.Range("H" + $newRow).Select() | Out-Null

# --- 5. Restore the original active sheet -----------------------------------
$wb.Worksheets.Item($originalActiveSheet).Activate() | Out-Null

$excel.CalculateFull() | Out-Null
